$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 12.31250333333333
$ws.Cells.Item(2, 8).Value = 36.93751
$ws.Cells.Item(2, 9).Value = 0.6498350963072504
$ws.Cells.Item(2, 10).Value = 0.6498350963072506
$ws.Cells.Item(2, 13).Value = 9.101967
$ws.Cells.Item(2, 14).Value = 27.305901
$ws.Cells.Item(2, 15).Value = 0.2872601673725235
$ws.Cells.Item(2, 16).Value = 0.2872601673725235
$ws.Cells.Item(2, 17).Value = 112.06799902739
$ws.Cells.Item(2, 18).Value = 1008.61199124651
$ws.Cells.Item(2, 19).Value = 0.1866717385297607
$ws.Cells.Item(2, 20).Value = 0.1866717385297607
$ws.Cells.Item(3, 7).Value = 12.31250333333333
$ws.Cells.Item(3, 8).Value = 36.93751
$ws.Cells.Item(3, 9).Value = 0.6498350963072504
$ws.Cells.Item(3, 10).Value = 0.6498350963072506
$ws.Cells.Item(3, 15).Value = 0.3055950511371977
$ws.Cells.Item(3, 16).Value = 0.3055950511371977
$ws.Cells.Item(3, 17).Value = 119.2209355263867
$ws.Cells.Item(3, 18).Value = 1072.98841973748
$ws.Cells.Item(3, 19).Value = 0.19858638948676
$ws.Cells.Item(3, 20).Value = 0.19858638948676
$ws.Cells.Item(4, 7).Value = 12.31250333333333
$ws.Cells.Item(4, 8).Value = 36.93751
$ws.Cells.Item(4, 9).Value = 0.6498350963072504
$ws.Cells.Item(4, 10).Value = 0.6498350963072506
$ws.Cells.Item(4, 13).Value = 3.905093666666666
$ws.Cells.Item(4, 14).Value = 11.715281
$ws.Cells.Item(4, 15).Value = 0.1232456523180152
$ws.Cells.Item(4, 16).Value = 0.1232456523180152
$ws.Cells.Item(4, 17).Value = 48.08147878781222
$ws.Cells.Item(4, 18).Value = 432.73330909031
$ws.Cells.Item(4, 19).Value = 0.08008935034352732
$ws.Cells.Item(4, 20).Value = 0.08008935034352734
$ws.Cells.Item(5, 7).Value = 12.31250333333333
$ws.Cells.Item(5, 8).Value = 36.93751
$ws.Cells.Item(5, 9).Value = 0.6498350963072504
$ws.Cells.Item(5, 10).Value = 0.6498350963072506
$ws.Cells.Item(5, 13).Value = 6.285238333333333
$ws.Cells.Item(5, 14).Value = 18.855715
$ws.Cells.Item(5, 15).Value = 0.1983635642284282
$ws.Cells.Item(5, 16).Value = 0.1983635642284282
$ws.Cells.Item(5, 17).Value = 77.3870179299611
$ws.Cells.Item(5, 18).Value = 696.48316136965
$ws.Cells.Item(5, 19).Value = 0.1289036058642301
$ws.Cells.Item(5, 20).Value = 0.1289036058642301
$ws.Cells.Item(6, 7).Value = 12.31250333333333
$ws.Cells.Item(6, 8).Value = 36.93751
$ws.Cells.Item(6, 9).Value = 0.6498350963072504
$ws.Cells.Item(6, 10).Value = 0.6498350963072506
$ws.Cells.Item(6, 13).Value = 2.710232666666667
$ws.Cells.Item(6, 14).Value = 8.130698000000001
$ws.Cells.Item(6, 15).Value = 0.08553556494383548
$ws.Cells.Item(6, 16).Value = 0.08553556494383548
$ws.Cells.Item(6, 17).Value = 33.36974874244223
$ws.Cells.Item(6, 18).Value = 300.3277386819801
$ws.Cells.Item(6, 19).Value = 0.0555840120829724
$ws.Cells.Item(6, 20).Value = 0.05558401208297241
$ws.Cells.Item(7, 9).Value = 0.3333514949915254
$ws.Cells.Item(7, 10).Value = 0.3333514949915254
$ws.Cells.Item(7, 13).Value = 9.101967
$ws.Cells.Item(7, 14).Value = 27.305901
$ws.Cells.Item(7, 15).Value = 0.2872601673725235
$ws.Cells.Item(7, 16).Value = 0.2872601673725235
$ws.Cells.Item(7, 17).Value = 57.488484738328
$ws.Cells.Item(7, 18).Value = 517.396362644952
$ws.Cells.Item(7, 19).Value = 0.0957586062451465
$ws.Cells.Item(7, 20).Value = 0.09575860624514651
$ws.Cells.Item(8, 9).Value = 0.3333514949915254
$ws.Cells.Item(8, 10).Value = 0.3333514949915254
$ws.Cells.Item(8, 15).Value = 0.3055950511371977
$ws.Cells.Item(8, 16).Value = 0.3055950511371977
$ws.Cells.Item(8, 19).Value = 0.1018705671585965
$ws.Cells.Item(8, 20).Value = 0.1018705671585965
$ws.Cells.Item(9, 9).Value = 0.3333514949915254
$ws.Cells.Item(9, 10).Value = 0.3333514949915254
$ws.Cells.Item(9, 13).Value = 3.905093666666666
$ws.Cells.Item(9, 14).Value = 11.715281
$ws.Cells.Item(9, 15).Value = 0.1232456523180152
$ws.Cells.Item(9, 16).Value = 0.1232456523180152
$ws.Cells.Item(9, 17).Value = 24.66476945674577
$ws.Cells.Item(9, 18).Value = 221.982925110712
$ws.Cells.Item(9, 19).Value = 0.04108412245141613
$ws.Cells.Item(9, 20).Value = 0.04108412245141613
$ws.Cells.Item(10, 9).Value = 0.3333514949915254
$ws.Cells.Item(10, 10).Value = 0.3333514949915254
$ws.Cells.Item(10, 13).Value = 6.285238333333333
$ws.Cells.Item(10, 14).Value = 18.855715
$ws.Cells.Item(10, 15).Value = 0.1983635642284282
$ws.Cells.Item(10, 16).Value = 0.1983635642284282
$ws.Cells.Item(10, 17).Value = 39.69788376540889
$ws.Cells.Item(10, 18).Value = 357.28095388868
$ws.Cells.Item(10, 19).Value = 0.06612479068739401
$ws.Cells.Item(10, 20).Value = 0.06612479068739402
$ws.Cells.Item(11, 9).Value = 0.3333514949915254
$ws.Cells.Item(11, 10).Value = 0.3333514949915254
$ws.Cells.Item(11, 13).Value = 2.710232666666667
$ws.Cells.Item(11, 14).Value = 8.130698000000001
$ws.Cells.Item(11, 15).Value = 0.08553556494383548
$ws.Cells.Item(11, 16).Value = 0.08553556494383548
$ws.Cells.Item(11, 17).Value = 17.11796684112178
$ws.Cells.Item(11, 18).Value = 154.061701570096
$ws.Cells.Item(11, 19).Value = 0.02851340844897227
$ws.Cells.Item(11, 20).Value = 0.02851340844897227
$ws.Cells.Item(12, 7).Value = 0.3185656666666667
$ws.Cells.Item(12, 8).Value = 0.955697
$ws.Cells.Item(12, 9).Value = 0.01681340870122405
$ws.Cells.Item(12, 10).Value = 0.01681340870122405
$ws.Cells.Item(12, 13).Value = 9.101967
$ws.Cells.Item(12, 14).Value = 27.305901
$ws.Cells.Item(12, 15).Value = 0.2872601673725235
$ws.Cells.Item(12, 16).Value = 0.2872601673725235
$ws.Cells.Item(12, 17).Value = 2.899574185333
$ws.Cells.Item(12, 18).Value = 26.096167667997
$ws.Cells.Item(12, 19).Value = 0.004829822597616263
$ws.Cells.Item(12, 20).Value = 0.004829822597616263
$ws.Cells.Item(13, 7).Value = 0.3185656666666667
$ws.Cells.Item(13, 8).Value = 0.955697
$ws.Cells.Item(13, 9).Value = 0.01681340870122405
$ws.Cells.Item(13, 10).Value = 0.01681340870122405
$ws.Cells.Item(13, 15).Value = 0.3055950511371977
$ws.Cells.Item(13, 16).Value = 0.3055950511371977
$ws.Cells.Item(13, 17).Value = 3.084644590817334
$ws.Cells.Item(13, 18).Value = 27.761801317356
$ws.Cells.Item(13, 19).Value = 0.005138094491841169
$ws.Cells.Item(13, 20).Value = 0.005138094491841169
$ws.Cells.Item(14, 7).Value = 0.3185656666666667
$ws.Cells.Item(14, 8).Value = 0.955697
$ws.Cells.Item(14, 9).Value = 0.01681340870122405
$ws.Cells.Item(14, 10).Value = 0.01681340870122405
$ws.Cells.Item(14, 13).Value = 3.905093666666666
$ws.Cells.Item(14, 14).Value = 11.715281
$ws.Cells.Item(14, 15).Value = 0.1232456523180152
$ws.Cells.Item(14, 16).Value = 0.1232456523180152
$ws.Cells.Item(14, 17).Value = 1.244028767317444
$ws.Cells.Item(14, 18).Value = 11.196258905857
$ws.Cells.Item(14, 19).Value = 0.002072179523071751
$ws.Cells.Item(14, 20).Value = 0.002072179523071751
$ws.Cells.Item(15, 7).Value = 0.3185656666666667
$ws.Cells.Item(15, 8).Value = 0.955697
$ws.Cells.Item(15, 9).Value = 0.01681340870122405
$ws.Cells.Item(15, 10).Value = 0.01681340870122405
$ws.Cells.Item(15, 13).Value = 6.285238333333333
$ws.Cells.Item(15, 14).Value = 18.855715
$ws.Cells.Item(15, 15).Value = 0.1983635642284282
$ws.Cells.Item(15, 16).Value = 0.1983635642284282
$ws.Cells.Item(15, 17).Value = 2.002261139817222
$ws.Cells.Item(15, 18).Value = 18.020350258355
$ws.Cells.Item(15, 19).Value = 0.00333516767680407
$ws.Cells.Item(15, 20).Value = 0.00333516767680407
$ws.Cells.Item(16, 7).Value = 0.3185656666666667
$ws.Cells.Item(16, 8).Value = 0.955697
$ws.Cells.Item(16, 9).Value = 0.01681340870122405
$ws.Cells.Item(16, 10).Value = 0.01681340870122405
$ws.Cells.Item(16, 13).Value = 2.710232666666667
$ws.Cells.Item(16, 14).Value = 8.130698000000001
$ws.Cells.Item(16, 15).Value = 0.08553556494383548
$ws.Cells.Item(16, 16).Value = 0.08553556494383548
$ws.Cells.Item(16, 17).Value = 0.8633870762784446
$ws.Cells.Item(16, 18).Value = 7.770483686506001
$ws.Cells.Item(16, 19).Value = 0.001438144411890798
$ws.Cells.Item(16, 20).Value = 0.001438144411890798
